$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# The department for this course row changed from the full faculty name
# to the short subject name.
$ws.Range("C2").Value = "English"

# The promotion-validity note has expired / no longer applies, so the
# cell is cleared out.
$ws.Range("R2").Value = ""

$ws.Range("R2").Select()
